# Apply cryptos list update (Sun Sep 10 11:54:48 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.965.81"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.634.51"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'214.18"
$ws.Range("E5").Value = "  -1.20%  "

$ws.Range("E6").Value = "  -0.75%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.253"
$ws.Range("E8").Value = "  -1.69%  "

$ws.Range("E9").Value = "  -2.96%  "

$ws.Range("D10").Value = "'18.48"
$ws.Range("E10").Value = "  -6.18%  "

$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").Value = "1.863.38"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.635.08"
$ws.Range("E13").Value = "  -0.61%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.18"
$ws.Range("E14").Value = "  -2.81%  "

$ws.Range("E15").Value = "  -2.64%  "

$ws.Range("D16").Value = "25.992.55"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("E17").Value = "  -2.95%  "

$ws.Range("D18").Value = "'61.65"
$ws.Range("E18").Value = "  -2.41%  "

$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").Value = "'190.30"
$ws.Range("E20").Value = "  -1.67%  "

$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("E22").Value = "  -3.84%  "

$ws.Range("E23").Value = "  -2.09%  "

$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").Value = "'143.37"
$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  -3.42%  "

$ws.Range("E28").Value = "  -1.97%  "

$ws.Range("E29").Value = "  -2.15%  "

$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  -1.34%  "

$ws.Range("D31").Value = "'0.0483"
$ws.Range("E31").Value = "  -3.54%  "

$ws.Range("E32").Value = "  -3.16%  "

$ws.Range("E33").Value = "  -4.58%  "

$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("E35").Value = "  -2.88%  "

$ws.Range("D36").Value = "1.138.87"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").Value = "'0.865"
$ws.Range("E37").Value = "  -4.55%  "

$ws.Range("D38").Value = "'2.44"
$ws.Range("E38").Value = "  -1.28%  "

$ws.Range("E39").Value = "  -3.54%  "

$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("D41").Value = "'98.49"
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").Value = "  -2.29%  "

$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.772.55"
$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.23"
$ws.Range("E44").Value = "  -5.02%  "

$ws.Range("E45").Value = "  -1.68%  "

$ws.Range("D46").Value = "'55.16"
$ws.Range("E46").Value = "  -2.97%  "

$ws.Range("E47").Value = "  -0.47%  "

$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").Value = "'7.55"
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("E51").Value = "  +0.09%  "
